$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7700
$ws.Range("I51").Value = 3500
$ws.Range("J51").Value = 10500
$ws.Range("K51").Value = 3500
$ws.Range("L51").Value = 10500
$ws.Range("M51").Value = -3016
$ws.Range("N51").Value = -11468

$ws.Range("H57").Value = 186566
$ws.Range("J57").Value = 186566
$ws.Range("L57").Value = 559698
$ws.Range("N57").Value = -560696

$ws.Range("H98").Value = 1952.12
$ws.Range("I98").Value = 1540.95
$ws.Range("K98").Value = 1540.95
$ws.Range("M98").Value = -42.95000000000005

$ws.Range("H100").Value = 4227.5625
$ws.Range("I100").Value = 3688.2222
$ws.Range("J100").Value = 4921
$ws.Range("K100").Value = 3688.2222
$ws.Range("L100").Value = 4921
$ws.Range("M100").Value = -3147.2222
$ws.Range("N100").Value = -6003

$ws.Range("H122").Value = 1952.12
$ws.Range("I122").Value = 1540.95
$ws.Range("K122").Value = 4622.85
$ws.Range("M122").Value = -2172.85

$ws.Range("H137").Value = 15768.083
$ws.Range("I137").Value = 9181.817999999999
$ws.Range("K137").Value = 27545.454
$ws.Range("M137").Value = -24995.454

$ws.Range("H140").Value = 78999.60000000001
$ws.Range("J140").Value = 78999.60000000001
$ws.Range("L140").Value = 78999.60000000001
$ws.Range("N140").Value = -89359.60000000001

$ws.Range("H141").Value = 3765.394
$ws.Range("I141").Value = 3891.9333
$ws.Range("K141").Value = 11675.7999
$ws.Range("M141").Value = -6495.7999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2253489
$ws.Range("I32").Value = 3345567.2
$ws.Range("K32").Value = 3345567.2
$ws.Range("M32").Value = -3345280.2

$ws.Range("H45").Value = 3573
$ws.Range("I45").Value = 1402.2
$ws.Range("K45").Value = 1402.2
$ws.Range("M45").Value = -1025.2

$ws.Range("H61").Value = 16187.111
$ws.Range("I61").Value = 4335.8
$ws.Range("J61").Value = 43121.91
$ws.Range("K61").Value = 4335.8
$ws.Range("L61").Value = 43121.91
$ws.Range("M61").Value = -4123.8
$ws.Range("N61").Value = -43545.91

$ws.Range("H74").Value = 20605.738
$ws.Range("I74").Value = 2749
$ws.Range("J74").Value = 32085.072
$ws.Range("K74").Value = 2749
$ws.Range("L74").Value = 32085.072
$ws.Range("M74").Value = -1875
$ws.Range("N74").Value = -33833.072

$ws.Range("H77").Value = 20605.738
$ws.Range("I77").Value = 2749
$ws.Range("J77").Value = 32085.072
$ws.Range("K77").Value = 13745
$ws.Range("L77").Value = 160425.36
$ws.Range("M77").Value = -9377
$ws.Range("N77").Value = -169161.36

$ws.Range("H132").Value = 6628.5625
$ws.Range("I132").Value = 2799.849
$ws.Range("K132").Value = 8399.547
$ws.Range("M132").Value = -5869.547

$ws.Range("H136").Value = 16187.111
$ws.Range("I136").Value = 4335.8
$ws.Range("J136").Value = 43121.91
$ws.Range("K136").Value = 13007.4
$ws.Range("L136").Value = 129365.73
$ws.Range("M136").Value = -10457.4
$ws.Range("N136").Value = -134465.73

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26074.744
$ws.Range("I31").Value = 17509.182
$ws.Range("K31").Value = 17509.182
$ws.Range("M31").Value = -17214.182

$ws.Range("H34").Value = 26074.744
$ws.Range("I34").Value = 17509.182
$ws.Range("K34").Value = 17509.182
$ws.Range("M34").Value = -17307.182

$ws.Range("H132").Value = 6920.1665
$ws.Range("I132").Value = 1666.85
$ws.Range("J132").Value = 13486.8125
$ws.Range("K132").Value = 5000.549999999999
$ws.Range("L132").Value = 40460.4375
$ws.Range("M132").Value = -2470.549999999999
$ws.Range("N132").Value = -45520.4375

$ws.Range("H134").Value = 5184.317
$ws.Range("I134").Value = 1732.2963
$ws.Range("K134").Value = 5196.8889
$ws.Range("M134").Value = -2661.8889

$ws.Range("H139").Value = 134389.5
$ws.Range("J139").Value = 200000
$ws.Range("L139").Value = 200000
$ws.Range("N139").Value = -210280

$ws.Range("H141").Value = 709916.75
$ws.Range("I141").Value = 99890.5
$ws.Range("J141").Value = 913258.8
$ws.Range("K141").Value = 99890.5
$ws.Range("L141").Value = 913258.8
$ws.Range("M141").Value = -94710.5
$ws.Range("N141").Value = -923618.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1472.0667
$ws.Range("I11").Value = 1164.8334
$ws.Range("J11").Value = 2701
$ws.Range("K11").Value = 3494.5002
$ws.Range("L11").Value = 8103
$ws.Range("M11").Value = -3354.5002
$ws.Range("N11").Value = -8383

$ws.Range("H12").Value = 39.857143
$ws.Range("I12").Value = 13.25
$ws.Range("K12").Value = 39.75
$ws.Range("M12").Value = 133.25

$ws.Range("H56").Value = 76929360
$ws.Range("I56").Value = 76929360
$ws.Range("K56").Value = 76929360
$ws.Range("M56").Value = -76928830

$ws.Range("H64").Value = 4686.385
$ws.Range("I64").Value = 3099.75
$ws.Range("K64").Value = 9299.25
$ws.Range("M64").Value = -9029.25

$ws.Range("H67").Value = 4686.385
$ws.Range("I67").Value = 3099.75
$ws.Range("K67").Value = 9299.25
$ws.Range("M67").Value = -8363.25

$ws.Range("H69").Value = 2250

$ws.Range("H72").Value = 2250

$ws.Range("H109").Value = 1967088.4
$ws.Range("J109").Value = 4775701.5
$ws.Range("L109").Value = 14327104.5
$ws.Range("N109").Value = -14329184.5

$ws.Range("H127").Value = 52718.43
$ws.Range("J127").Value = 52718.43
$ws.Range("L127").Value = 158155.29
$ws.Range("N127").Value = -168075.29

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 13055.143
$ws.Range("I126").Value = 3699.5
$ws.Range("K126").Value = 11098.5
$ws.Range("M126").Value = -8628.5

$ws.Range("H132").Value = 9692.23
$ws.Range("I132").Value = 12466.667
$ws.Range("J132").Value = 3449.75
$ws.Range("K132").Value = 37400.001
$ws.Range("L132").Value = 10349.25
$ws.Range("M132").Value = -34870.001
$ws.Range("N132").Value = -15409.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6911.1562
$ws.Range("I40").Value = 3791.5
$ws.Range("K40").Value = 3791.5
$ws.Range("M40").Value = -3655.5

$ws.Range("H61").Value = 3490.7273
$ws.Range("I61").Value = 1730.6
$ws.Range("J61").Value = 4957.5
$ws.Range("K61").Value = 1730.6
$ws.Range("L61").Value = 4957.5
$ws.Range("M61").Value = -1528.6
$ws.Range("N61").Value = -5361.5

$ws.Range("H80").Value = 16750
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 16750
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 16750
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -18996

$ws.Range("H83").Value = 16750
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 16750
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 50250
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -61482

$ws.Range("H113").Value = 3490.7273
$ws.Range("I113").Value = 1730.6
$ws.Range("J113").Value = 4957.5
$ws.Range("K113").Value = 1730.6
$ws.Range("L113").Value = 4957.5
$ws.Range("M113").Value = 439.4000000000001
$ws.Range("N113").Value = -9297.5

$ws.Range("H132").Value = 11977.934
$ws.Range("I132").Value = 5190
$ws.Range("K132").Value = 15570
$ws.Range("M132").Value = -13040

$ws.Range("H136").Value = 20263.062
$ws.Range("I136").Value = 16329.5
$ws.Range("J136").Value = 26819
$ws.Range("K136").Value = 48988.5
$ws.Range("L136").Value = 80457
$ws.Range("M136").Value = -46438.5
$ws.Range("N136").Value = -85557

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 232185.6
$ws.Range("J138").Value = 232185.6
$ws.Range("L138").Value = 232185.6
$ws.Range("N138").Value = -242465.6

$ws.Range("H139").Value = 97088.664
$ws.Range("J139").Value = 97088.664
$ws.Range("L139").Value = 97088.664
$ws.Range("N139").Value = -107368.664

$ws.Range("H140").Value = 171275.1
$ws.Range("J140").Value = 180402.6
$ws.Range("L140").Value = 180402.6
$ws.Range("N140").Value = -190762.6

$ws.Range("H141").Value = 83999.664
$ws.Range("J141").Value = 83999.664
$ws.Range("L141").Value = 83999.664
$ws.Range("N141").Value = -94359.664
